$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new "Sheet2" after "Sheet1"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with the lower/upper age-bucket table
$ws2.Range("A1").Value = "lower"
$ws2.Range("B1").Value = "upper"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 19

$ws2.Range("A3").Value = 20
$ws2.Range("B3").Value = 29

$ws2.Range("A4").Value = 30
$ws2.Range("B4").Value = 39

$ws2.Range("A5").Value = 40
$ws2.Range("B5").Value = 49

$ws2.Range("A6").Value = 50
$ws2.Range("B6").Value = 59

$ws2.Range("A7").Value = 60
$ws2.Range("B7").Value = 100

# Widen column D on Sheet1 to fit its (now visible/used) contents
$ws1.Columns.Item(4).AutoFit()

# Restore focus/selection on Sheet1 before Sheet2 becomes the active sheet
$ws1.Activate()
[void]$ws1.Range("B105").Select()

# Make Sheet2 the active tab and match its final selection
$ws2.Activate()
[void]$ws2.Range("G37").Select()
